$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# The sheet's Hyperlinks collection addresses its entries by row number.
# We're about to insert a row (shifting today's rows 15-20 down to 16-21),
# so drop every existing hyperlink object now and recreate them once all
# the URLs are sitting in their final rows/cells.
$ws.Hyperlinks.Delete()

# Insert one blank row at row 15; Excel pushes the existing rows 15-20
# (values + styles) down to rows 16-21.
$ws.Rows.Item(15).Insert()

$newTimestamp = "2025-12-16 01:57:15"

# Refresh the "取得日時" (fetched-at) column for every data row (2-21) to
# match the new scrape time recorded in this commit.
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Populate the freshly inserted row 15 with the new listing that was
# scraped in this run.
$ws.Cells.Item(15, 2).Value = "【急募】帳票デジタル化のフロントエンド開発者募集"
$ws.Cells.Item(15, 3).Value = "システム開発"
$ws.Cells.Item(15, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(15, 5).Value = "期限情報なし"
$ws.Cells.Item(15, 6).Value = "https://www.lancers.jp/work/detail/5454857"
$ws.Cells.Item(15, 7).Value = 68
$ws.Cells.Item(15, 8).Value = "◆開発"

# Recreate the column-F hyperlinks for every data row now that all the
# target URLs live in their final positions (this also re-applies the
# hyperlink font styling used throughout column F).
for ($r = 2; $r -le 21; $r++) {
    $target = $ws.Cells.Item($r, 6).Value()
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $target)
}

Write-Output "Edit complete"
